$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated title timestamp
$ws.Range("A1").Value = 'Datos actualizados a 25 de Marzo de 2020 a las 11:46'

# Full refreshed dataset (re-sorted by total cases, columns B:H updated)
$data = @(
    @('China',81218,47,73650,4287,1399,4,3281),
    @('Italia',69176,0,8326,54030,3393,0,6820),
    @('Estados Unidos',54941,60,379,53778,1175,4,784),
    @('España',47610,5552,5367,38809,2636,443,3434),
    @('Alemania',34009,1018,3532,30305,23,13,172),
    @('Iran',27017,2206,9625,15315,0,143,2077),
    @('Francia',22304,0,3281,17923,2516,0,1100),
    @('Suiza',10171,294,131,9905,141,13,135),
    @('Corea del Sur',9137,100,3730,5281,59,6,126),
    @('Reino Unido',8077,0,135,7520,20,0,422),
    @('Paises Bajos',5560,0,2,5282,435,0,276),
    @('Austria',5494,211,9,5455,26,2,30),
    @('Belgica',4937,668,547,4212,474,56,178),
    @('Noruega',2868,2,6,2849,44,1,13),
    @('Canada',2792,0,112,2654,1,0,26),
    @('Australia',2423,106,118,2296,11,1,9),
    @('Portugal',2362,0,22,2307,48,0,33),
    @('Suecia',2300,1,16,2244,144,0,40),
    @('Brasil',2271,24,2,2222,18,1,47),
    @('Israel',2030,100,58,1967,37,2,5),
    @('Turquia',1872,0,0,1828,0,0,44),
    @('Malasia',1796,172,199,1578,64,3,19),
    @('Dinamarca',1715,124,1,1682,69,0,32),
    @('Chequia',1497,103,10,1484,19,0,3),
    @('Irlanda',1329,0,5,1317,29,0,7),
    @('Japon',1193,0,285,865,54,0,43),
    @('Luxemburgo',1099,0,6,1085,3,0,8),
    @('Ecuador',1082,0,3,1052,2,0,27),
    @('Pakistan',1000,28,19,974,5,0,7),
    @('Tailandia',934,107,70,860,11,0,4),
    @('Polonia',927,26,1,914,3,2,12),
    @('Chile',922,0,17,903,7,0,2),
    @('Finlandia',841,49,10,828,11,2,3),
    @('Rumania',794,0,79,702,15,1,13),
    @('Indonesia',790,104,31,701,0,3,58),
    @('Arabia Saudita',767,0,28,738,0,0,1),
    @('Grecia',743,0,29,694,35,0,20),
    @('Crucero',712,0,587,115,15,0,10),
    @('Sudafrica',709,155,4,705,2,0,0),
    @('Rusia',658,163,29,628,8,0,1),
    @('Islandia',648,0,51,595,13,0,2),
    @('Filipinas',636,84,26,572,1,3,38),
    @('India',562,26,40,512,0,0,10),
    @('Singapur',558,0,156,400,17,0,2),
    @('Catar',526,0,41,485,6,0,0),
    @('Eslovenia',480,0,3,473,12,0,4),
    @('Panama',443,0,1,436,33,0,6),
    @('Croacia',418,36,16,401,6,0,1),
    @('Peru',416,0,1,408,9,0,7),
    @('Hong Kong',410,23,102,304,4,0,4),
    @('Mexico',405,38,4,396,1,1,5),
    @('Estonia',404,35,8,396,5,0,0),
    @('Egipto',402,0,80,302,0,0,20),
    @('Barein',392,0,177,212,2,0,3),
    @('Argentina',387,0,52,329,0,0,6),
    @('Colombia',378,0,6,369,0,0,3),
    @('Libano',333,15,8,321,4,0,4),
    @('Irak',316,0,75,214,0,0,27),
    @('Republica Dominicana',312,0,3,303,0,0,6),
    @('Serbia',303,0,15,284,21,1,4),
    @('Armenia',265,16,16,249,6,0,0),
    @('Argelia',264,0,65,180,0,0,19),
    @('Lituania',255,46,1,250,1,2,4),
    @('Emiratos Arabes Unidos',248,0,45,201,2,0,2),
    @('Taiwan',235,19,29,204,0,0,2),
    @('Hungria',226,39,21,195,6,1,10),
    @('Letonia',221,24,1,220,0,0,0),
    @('Bulgaria',220,2,4,213,8,0,3),
    @('Eslovaquia',216,12,7,209,2,0,0),
    @('Nueva Zelanda',205,0,22,183,0,0,0),
    @('Kuwait',195,4,43,152,6,0,0),
    @('Uruguay',189,0,0,189,3,0,0),
    @('San Marino',187,0,4,162,12,0,21),
    @('Costa Rica',177,0,2,173,4,0,2),
    @('Marruecos',170,0,6,159,1,0,5),
    @('Bosnia y Herzegovina',168,0,2,163,1,0,3),
    @('Principado de Andorra',164,0,1,162,7,0,1),
    @('Jordania',154,0,1,153,0,0,0),
    @('Republica de Macedonia',148,0,1,145,1,0,2),
    @('Albania',146,23,17,124,3,0,5),
    @('Vietnam',134,0,17,117,3,0,0),
    @('Islas Feroe',132,10,38,94,2,0,0),
    @('Moldavia',125,0,2,122,20,0,1),
    @('Republica de Chipre',124,0,3,118,3,0,3),
    @('Tunez',114,0,1,109,11,0,4),
    @('Burkina Faso',114,0,7,103,0,0,4),
    @('Ucrania',113,11,1,108,0,1,4),
    @('Malta',110,0,2,108,1,0,0),
    @('Brunei',109,5,2,107,1,0,0),
    @('Sri Lanka',102,0,3,99,2,0,0),
    @('Senegal',99,13,9,90,0,0,0),
    @('Oman',99,15,17,82,0,0,0),
    @('Reunion',94,0,1,93,0,0,0),
    @('Camboya',93,2,4,89,1,0,0),
    @('Venezuela',91,7,15,76,2,0,0),
    @('Azerbaiyan',87,0,10,76,6,0,1),
    @('Bielorrusia',86,5,29,57,2,0,0),
    @('Kazajistan',79,7,0,79,0,0,0),
    @('Afganistan',74,0,1,71,0,1,2),
    @('Guadalupe',73,0,0,72,4,0,1),
    @('Costa de Marfil',73,0,2,71,0,0,0),
    @('Georgia',73,3,10,63,1,0,0),
    @('Camerun',66,0,2,63,0,0,1),
    @('Estado de Palestina',62,2,16,46,0,0,0),
    @('Trinidad yTobago',57,0,0,57,0,0,0),
    @('Martinica',57,0,0,56,7,0,1),
    @('Uzbekistan',55,5,0,55,4,0,0),
    @('Ghana',53,0,0,51,0,0,2),
    @('Montenegro',52,5,0,51,0,0,1),
    @('Liechtenstein',51,0,0,51,0,0,0),
    @('Mauricio',48,6,0,46,1,0,2),
    @('Cuba',48,0,1,46,2,0,1),
    @('Nigeria',46,2,2,43,0,0,1),
    @('Consejo Danes para los Refugiados',45,0,0,43,0,0,2),
    @('Kirguistan',44,2,0,44,0,0,0),
    @('Ruanda',40,0,0,40,0,0,0),
    @('Puerto Rico',39,0,1,36,0,0,2),
    @('Banglades',39,0,7,27,0,1,5),
    @('Paraguay',37,10,0,34,1,1,3),
    @('Mayotte',36,0,0,36,0,0,0),
    @('Honduras',36,6,0,36,0,0,0),
    @('Bolivia',32,3,0,32,0,0,0),
    @('Guam',32,0,0,31,0,0,1),
    @('Macao',30,1,10,20,0,0,0),
    @('Polinesia Francesa',25,0,0,25,0,0,0),
    @('Kenia',25,0,0,25,0,0,0),
    @('Jamaica',25,4,2,22,0,0,1),
    @('Isla de Man',23,0,0,23,0,0,0),
    @('Monaco',23,0,1,22,0,0,0),
    @('Guayana Francesa',23,0,6,17,0,0,0),
    @('Guatemala',21,0,0,20,0,0,1),
    @('Togo',20,0,1,19,0,0,0),
    @('Madagascar',19,2,0,19,0,0,0),
    @('Barbados',18,0,0,18,0,0,0),
    @('Islas Virgenes de los Estados Unidos',17,0,0,17,0,0,0),
    @('Aruba',17,0,1,16,0,0,0),
    @('Gibraltar',15,0,5,10,0,0,0),
    @('Nueva Caledonia',14,4,0,14,0,0,0),
    @('Maldivas',13,0,5,8,0,0,0),
    @('Tanzania',12,0,0,12,0,0,0),
    @('Zambia',12,9,0,12,0,0,0),
    @('Etiopia',12,0,0,12,0,0,0),
    @('Mongolia',10,0,0,10,0,0,0),
    @('Uganda',9,0,0,9,0,0,0),
    @('Guinea Ecuatorial',9,0,0,9,0,0,0),
    @('El Salvador',9,4,0,9,0,0,0),
    @('San Martin (Parte Francesa)',8,0,0,8,0,0,0),
    @('Haiti',7,0,0,7,0,0,0),
    @('Dominica',7,0,0,7,0,0,0),
    @('Surinam',7,0,0,7,0,0,0),
    @('Seychelles',7,0,0,7,0,0,0),
    @('Niger',7,4,0,6,0,1,1),
    @('Namibia',7,0,2,5,0,0,0),
    @('Benin',6,0,0,6,0,0,0),
    @('Bermudas',6,0,0,6,0,0,0),
    @('Islas Caimanes',6,0,0,5,0,0,1),
    @('Gabon',6,0,0,5,0,0,1),
    @('Curazao',6,0,0,5,0,0,1),
    @('Fiyi',5,1,0,5,0,0,0),
    @('Bahamas',5,0,1,4,0,0,0),
    @('Guyana',5,0,0,4,0,0,1),
    @('Groenlandia',5,0,2,3,0,0,0),
    @('Congo',4,0,0,4,0,0,0),
    @('Guinea',4,0,0,4,0,0,0),
    @('Santa Sede',4,0,0,4,0,0,0),
    @('Suazilandia',4,0,0,4,0,0,0),
    @('Birmania',3,0,0,3,0,0,0),
    @('Angola',3,0,0,3,0,0,0),
    @('Liberia',3,0,0,3,0,0,0),
    @('Santa Lucia',3,0,0,3,0,0,0),
    @('Antigua y Barbuda',3,0,0,3,0,0,0),
    @('San Bartolome',3,0,0,3,0,0,0),
    @('Mozambique',3,0,0,3,0,0,0),
    @('Republica de Africa Central',3,0,0,3,0,0,0),
    @('Republica de Yibuti',3,0,0,3,0,0,0),
    @('Republica del Chad',3,0,0,3,0,0,0),
    @('Nepal',3,1,1,2,0,0,0),
    @('Zimbabue',3,0,0,2,0,0,1),
    @('Sudan',3,0,0,2,0,0,1),
    @('Gambia',3,0,0,2,0,0,1),
    @('Cabo Verde',3,0,0,2,0,0,1),
    @('Mauritania',2,0,0,2,0,0,0),
    @('Butan',2,0,0,2,0,0,0),
    @('Laos',2,0,0,2,0,0,0),
    @('Nicaragua',2,0,0,2,0,0,0),
    @('San Martin (Parte Holandesa)',2,0,0,2,0,0,0),
    @('Mali',2,2,0,2,0,0,0),
    @('Siria',1,0,0,1,0,0,0),
    @('Libia',1,0,0,1,0,0,0),
    @('San Vicente y las Granadinas',1,0,0,1,0,0,0),
    @('Granada',1,0,0,1,0,0,0),
    @('Montserrat',1,0,0,1,0,0,0),
    @('Somalia',1,0,0,1,0,0,0),
    @('Papua Nueva Guinea',1,0,0,1,0,0,0),
    @('Timor Oriental',1,0,0,1,0,0,0),
    @('Eritrea',1,0,0,1,0,0,0),
    @('Islas Turcas y Caicos',1,0,0,1,0,0,0),
    @('Belice',1,0,0,1,0,0,0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
}
